$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTD")

# Fill in the newly-completed test case row (row 19 / ID 18)
$ws.Range("B19").Value = 10
$ws.Range("C19").Value = "All fields detect user input (click)"
$ws.Range("D19").Value = "actionlisteners + html elements"
$ws.Range("E19").Value = 'console.log("HELLO :)");'
$ws.Range("F19").Value = "Add all eventlisteners to all inputs"
$ws.Range("G19").Value = "hello in console for each input clicked on"
$ws.Range("H19").Value = "All fields produce hello when clicked "
$ws.Range("I19").Value = "pass"
$ws.Range("J19").Value = "Changed array to an object with help from https://masteringjs.io/tutorials/fundamentals/foreach-object to add event listeners, while keeping code easier to read"
$ws.Range("K19").Value = "https://imgur.com/EDmnBzc"

# Row now holds wrapped multi-line text, same visual height as similar rows above
$ws.Rows.Item(19).RowHeight = 45

# Turn the proof screenshot link into a real hyperlink (like the rows above it)
$ws.Hyperlinks.Add($ws.Range("K19"), "https://imgur.com/EDmnBzc")

# Adding the hyperlink re-styles the cell to the workbook's "Hyperlink" look;
# restore the plain style used by the rest of the K column (copy format only)
$ws.Range("K18").Copy()
$ws.Range("K19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to where the user last clicked
$ws.Range("M16").Select() | Out-Null
